$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = -12.4815
$ws.Range("A10").Value = -21.65329999999999
$ws.Range("A12").Value = -21.57940000000001
$ws.Range("C15").Value = -14.69949999999998
$ws.Range("A18").Value = -22.16750000000001
$ws.Range("C20").Value = -12.1928
$ws.Range("C29").Value = -11.1562
$ws.Range("C30").Value = -12.92649999999999
$ws.Range("C31").Value = -12.9578
$ws.Range("A37").Value = -20.10750000000001
$ws.Range("C40").Value = -13.05380000000001
$ws.Range("A55").Value = -22.2729
$ws.Range("A68").Value = -21.76199999999999
$ws.Range("C68").Value = -11.8639
$ws.Range("C76").Value = -12.52630000000001
$ws.Range("A77").Value = -20.71819999999999
$ws.Range("A78").Value = -20.29939999999998
$ws.Range("C87").Value = -13.59059999999999
$ws.Range("C88").Value = -12.61839999999999
$ws.Range("C96").Value = -12.96500000000001
$ws.Range("C98").Value = -12.5369
$ws.Range("C101").Value = -12.68500000000001
$ws.Range("C102").Value = -13.15400000000001
